$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A7:P16").EntireRow.Delete()
Write-Host "done"
